$d = $word.ActiveDocument

# Find the table that contains the "otherParcels" merge fields (the
# "Other parcels" detail table), then remove its "Certificate Of Title"
# row (label cell "Certificate Of Title" / value cell referencing
# {d.otherParcels[i].certificateOfTitle:...}).
for ($ti = $d.Tables.Count; $ti -ge 1; $ti--) {
    $t = $d.Tables.Item($ti)
    for ($ri = $t.Rows.Count; $ri -ge 1; $ri--) {
        $row = $t.Rows.Item($ri)
        $label = $row.Cells.Item(1).Range.Text
        $value = $row.Cells.Item(2).Range.Text
        if ($label -like "*Certificate Of Title*" -and $value -like "*otherParcels*certificateOfTitle*") {
            $row.Delete()
        }
    }
}
